# ==========================================================================
# Edit script: Fill in Social Engineering Attack risk table (table 1)
# ==========================================================================
$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-CellHighlightXml($cell, $innerXml) {
    $cell.Range.InsertXML($innerXml)
}

$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------
# Row 2 (already had content) - highlight everything yellow, tweak text
# ---------------------------------------------------------------------

# Cell 1: insert "or " before "golden ticket attack" (creates 2nd sentence piece)
$null = $d.Content.Find.Execute(
    "pass-the-hash attack, golden ticket attack, or a KDC attack.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "pass-the-hash attack, or golden ticket attack, or a KDC attack.", 2)

# Rebuild cell 1 paragraph with 3 runs + highlight (matches diff run split)
$cell = $t.Cell(2,1)
$xml = "<w:p $W>" +
  "<w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`">An attacker profiled a user on social media, crafted a phishing email and a custom payload. The phishing email contains a payload which was executed by the user. The attack was able to gain access to internal active directory using pass-the-hash attack, </w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`">or </w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>golden ticket attack, or a KDC attack.</w:t></w:r>" +
  "</w:p>"
Set-CellHighlightXml $cell $xml

# Cell 2: Moderate
$cell = $t.Cell(2,2)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Moderate</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

# Cell 3: Major
$cell = $t.Cell(2,3)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Major</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

# Cell 4: Medium
$cell = $t.Cell(2,4)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Medium</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

# Cell 5: Controls bullets (each bullet text + separate "." run), highlighted
$cell = $t.Cell(2,5)
$items = @(
  "Security awareness training",
  "Restrict user privileges on their workstation",
  "Set group policy for account lockout",
  "Harden Active Directory domain controllers"
)
$xml = ""
foreach ($it in $items) {
    $xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
      "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>$it</w:t></w:r>" +
      "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>.</w:t></w:r></w:p>"
}
Set-CellHighlightXml $cell $xml

Write-Host "Row 2 done"

# ---------------------------------------------------------------------
# Row 3 (new): USB scenario
# ---------------------------------------------------------------------

$cell = $t.Cell(3,1)
$xml = "<w:p $W>" +
  "<w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`">An attacker prepared a USB with malicious content. This USB was mailed to the user with instructions to plug-in, etc. Without knowing what was </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>actually installed</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> in the USB,</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>user</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> plugged the USB into the company’s network. The payload was then executed, and the attack started and was able to access the internal active directory using pass-the-hash attack, or golden ticket attack, or KDC attack.</w:t></w:r>" +
  "</w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(3,2)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Rare</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(3,3)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Major</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(3,4)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Medium</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(3,5)
$xml = ""
# 1. Security awareness training.
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Security awareness training</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>.</w:t></w:r></w:p>"
# 2. Think before you plug in a device.
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Think before you plug in a device</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>.</w:t></w:r></w:p>"
# 3. Seek advise from team members if the action is safe.
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`">Seek </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>advise</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/><w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> from team members</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> if the action is safe</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>.</w:t></w:r></w:p>"
# 4. Install tools that can detect malicious software trying to read passwords from memory.
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`">Install tools that can detect malicious </w:t></w:r>" +
  "<w:proofErr w:type=`"gramStart`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>software</w:t></w:r>" +
  "<w:proofErr w:type=`"gramEnd`"/>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> trying to read passwords from memory</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>.</w:t></w:r></w:p>"
# 5. Limit who can run admin tools.
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Limit who can run admin tools</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>.</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

Write-Host "Row 3 done"

# ---------------------------------------------------------------------
# Row 4 (new): IT / VPN scenario
# ---------------------------------------------------------------------

$cell = $t.Cell(4,1)
$xml = "<w:p $W>" +
  "<w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`">An attacker pretended to be from IT department and had a phone call to the user. Then the user unintentionally provided the necessary credentials </w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>for</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> the attacker to access the VPN on the network. </w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>With this action, the attacker was able to access the internal active directory using KDC attack, golden ticket attack, and pass-the-hash attack.</w:t></w:r>" +
  "</w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(4,2)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Likely</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(4,3)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Major</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(4,4)
$xml = "<w:p $W><w:pPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>High</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

$cell = $t.Cell(4,5)
$xml = ""
# 1. Security awareness training. (single run, period included)
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Security awareness training.</w:t></w:r></w:p>"
# 2. Restrict user privileges on their workstation.
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Restrict user privileges on their workstation.</w:t></w:r></w:p>"
# 3. Implement Multi-Factor Authentication (MFA) when accessing VPN.
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Implement Multi-Factor Authentication (MFA) when accessing VPN.</w:t></w:r></w:p>"
# 4. Choose a reputable VPN provider that uses modern encryption standard. (2 runs)
$xml += "<w:p $W><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"8`"/></w:numPr><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t>Choose a reputable VPN provider</w:t></w:r>" +
  "<w:r><w:rPr><w:highlight w:val=`"yellow`"/></w:rPr><w:t xml:space=`"preserve`"> that uses modern encryption standard.</w:t></w:r></w:p>"
Set-CellHighlightXml $cell $xml

Write-Host "Row 4 done"
